$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
try {
  $r = $tr.InsertAfter("X")
  Write-Output ("insertAfter ok: " + $tr.Text)
} catch { Write-Output ("ERR InsertAfter: " + $_) }
